# Append second batch of students to the Student_Emails sheet (rows 66-124)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 66
$rowCount = 59

# Columns: A=No. B=Student Number C=Student Name D=DoB F=Gender G=email (E left blank)
$data = New-Object 'object[,]' $rowCount,7

$data[0,0] = 1
$data[0,1] = 136603
$data[0,2] = 'Adhiambo, Jane Daisy'
$data[0,3] = 37684
$data[0,5] = 'F'
$data[0,6] = 'adaisy@gmail.com'

$data[1,0] = 2
$data[1,1] = 142285
$data[1,2] = 'Ali, Sumeiya Abdulle Abdirizak'
$data[1,3] = 37377
$data[1,5] = 'F'
$data[1,6] = 'aabdirizak@gmail.com'

$data[2,0] = 3
$data[2,1] = 146072
$data[2,2] = 'Asava, Wayne Majani'
$data[2,3] = 37563
$data[2,5] = 'M'
$data[2,6] = 'amajani@gmail.com'

$data[3,0] = 4
$data[3,1] = 145210
$data[3,2] = 'Barasa, Michelle Nekesa'
$data[3,3] = 36702
$data[3,5] = 'F'
$data[3,6] = 'bnekesa@gmail.com'

$data[4,0] = 5
$data[4,1] = 143496
$data[4,2] = 'Bigingi, Ian Duncan'
$data[4,3] = 37730
$data[4,5] = 'M'
$data[4,6] = 'bduncan@gmail.com'

$data[5,0] = 6
$data[5,1] = 146171
$data[5,2] = 'Cheboi, Lormotum Noel'
$data[5,3] = 36728
$data[5,5] = 'M'
$data[5,6] = 'cnoel@gmail.com'

$data[6,0] = 7
$data[6,1] = 145835
$data[6,2] = 'Dulo, Humphrey James'
$data[6,3] = 37368
$data[6,5] = 'M'
$data[6,6] = 'djames@gmail.com'

$data[7,0] = 8
$data[7,1] = 146565
$data[7,2] = 'Etemesi, Caleb Asira'
$data[7,3] = 37656
$data[7,5] = 'M'
$data[7,6] = 'easira@gmail.com'

$data[8,0] = 9
$data[8,1] = 136046
$data[8,2] = 'Frances, Ayango'
$data[8,3] = 37220
$data[8,5] = 'F'
$data[8,6] = 'fayango@gmail.com'

$data[9,0] = 10
$data[9,1] = 146399
$data[9,2] = 'Gitau, Sandra Wanjiru'
$data[9,3] = 37632
$data[9,5] = 'F'
$data[9,6] = 'gwanjiru@gmail.com'

$data[10,0] = 11
$data[10,1] = 146424
$data[10,2] = 'Hassan, Masoud Ali'
$data[10,3] = 37156
$data[10,5] = 'M'
$data[10,6] = 'hali1@gmail.com'

$data[11,0] = 12
$data[11,1] = 145279
$data[11,2] = 'Kabui, Michelle Miceere'
$data[11,3] = 37314
$data[11,5] = 'F'
$data[11,6] = 'kmiceere@gmail.com'

$data[12,0] = 13
$data[12,1] = 145491
$data[12,2] = 'Kamau, Edinah Nyambura'
$data[12,3] = 37660
$data[12,5] = 'F'
$data[12,6] = 'knyambura@gmail.com'

$data[13,0] = 14
$data[13,1] = 137991
$data[13,2] = 'Kamau, Jesse Mbugua'
$data[13,3] = 37847
$data[13,5] = 'M'
$data[13,6] = 'kmbugua@gmail.com'

$data[14,0] = 15
$data[14,1] = 146013
$data[14,2] = 'Karani, Amanda Ngendo'
$data[14,3] = 37613
$data[14,5] = 'F'
$data[14,6] = 'kngendo@gmail.com'

$data[15,0] = 16
$data[15,1] = 145212
$data[15,2] = 'Kathurima, Ryan Kinoti'
$data[15,3] = 36644
$data[15,5] = 'M'
$data[15,6] = 'kkinoti@gmail.com'

$data[16,0] = 17
$data[16,1] = 145802
$data[16,2] = 'Kennedy, Angel Venoliah'
$data[16,3] = 37250
$data[16,5] = 'F'
$data[16,6] = 'kvenoliah@gmail.com'

$data[17,0] = 18
$data[17,1] = 136395
$data[17,2] = 'Kiilu, Jeff Kioko'
$data[17,3] = 37618
$data[17,5] = 'M'
$data[17,6] = 'kkioko@gmail.com'

$data[18,0] = 19
$data[18,1] = 137503
$data[18,2] = 'Kiiru, Cynthia Everlyn Muthoni'
$data[18,3] = 37230
$data[18,5] = 'F'
$data[18,6] = 'kmuthoni1@gmail.com'

$data[19,0] = 20
$data[19,1] = 141690
$data[19,2] = 'Kimani, Alex Mwangi'
$data[19,3] = 37425
$data[19,5] = 'M'
$data[19,6] = 'kmwangi@gmail.com'

$data[20,0] = 21
$data[20,1] = 128576
$data[20,2] = 'Kinegeni, Terry'
$data[20,3] = 36869
$data[20,5] = 'F'
$data[20,6] = 'kterry@gmail.com'

$data[21,0] = 22
$data[21,1] = 145602
$data[21,2] = 'Kyalo, Felicia Mutheu'
$data[21,3] = 36953
$data[21,5] = 'F'
$data[21,6] = 'kmutheu@gmail.com'

$data[22,0] = 23
$data[22,1] = 138216
$data[22,2] = 'Leting, Sylvester Kiplagat'
$data[22,3] = 37970
$data[22,5] = 'M'
$data[22,6] = 'lkiplagat@gmail.com'

$data[23,0] = 24
$data[23,1] = 139991
$data[23,2] = 'Lihanda, Glen Musa'
$data[23,3] = 37951
$data[23,5] = 'M'
$data[23,6] = 'lmusa@gmail.com'

$data[24,0] = 25
$data[24,1] = 144915
$data[24,2] = 'Mahia, Jerome Kamau'
$data[24,3] = 37951
$data[24,5] = 'M'
$data[24,6] = 'mkamau@gmail.com'

$data[25,0] = 26
$data[25,1] = 115104
$data[25,2] = 'Mbugua, Nathan Ng''ethe'
$data[25,3] = 37433
$data[25,5] = 'M'
$data[25,6] = 'mngethe@gmail.com'

$data[26,0] = 27
$data[26,1] = 139074
$data[26,2] = 'Mbwanga, Emmanuel Chivunira'
$data[26,3] = 37567
$data[26,5] = 'M'
$data[26,6] = 'mchivunira@gmail.com'

$data[27,0] = 28
$data[27,1] = 92313
$data[27,2] = 'Moire, Henry Nyakundi'
$data[27,3] = 37718
$data[27,5] = 'M'
$data[27,6] = 'mnyakundi@gmail.com'

$data[28,0] = 29
$data[28,1] = 129029
$data[28,2] = 'Morara, Keith Matwere'
$data[28,3] = 36764
$data[28,5] = 'M'
$data[28,6] = 'mmatwere@gmail.com'

$data[29,0] = 30
$data[29,1] = 138583
$data[29,2] = 'Muchiri, Lynn Wairimu'
$data[29,3] = 37094
$data[29,5] = 'F'
$data[29,6] = 'mwairimu@gmail.com'

$data[30,0] = 31
$data[30,1] = 144338
$data[30,2] = 'Mugendi, Emmanuel Muthomi'
$data[30,3] = 36754
$data[30,5] = 'M'
$data[30,6] = 'mmuthomi@gmail.com'

$data[31,0] = 32
$data[31,1] = 140091
$data[31,2] = 'Mukiri, Maryanne Wanjiku'
$data[31,3] = 36812
$data[31,5] = 'F'
$data[31,6] = 'mwanjiku@gmail.com'

$data[32,0] = 33
$data[32,1] = 145351
$data[32,2] = 'Mungai, Kihanya'
$data[32,3] = 36779
$data[32,5] = 'M'
$data[32,6] = 'mkihanya@gmail.com'

$data[33,0] = 34
$data[33,1] = 145836
$data[33,2] = 'Musyoka, Brian Kioko'
$data[33,3] = 37447
$data[33,5] = 'M'
$data[33,6] = 'mkioko@gmail.com'

$data[34,0] = 35
$data[34,1] = 139133
$data[34,2] = 'Mutende, Arabella Fanisheba'
$data[34,3] = 37382
$data[34,5] = 'F'
$data[34,6] = 'mfanisheba@gmail.com'

$data[35,0] = 36
$data[35,1] = 145703
$data[35,2] = 'Mutinda, Bryan Lwaya'
$data[35,3] = 37773
$data[35,5] = 'M'
$data[35,6] = 'mlwaya@gmail.com'

$data[36,0] = 37
$data[36,1] = 146016
$data[36,2] = 'Mwago, Megan Dette'
$data[36,3] = 37117
$data[36,5] = 'F'
$data[36,6] = 'mdette@gmail.com'

$data[37,0] = 38
$data[37,1] = 145041
$data[37,2] = 'Mwai, David King'
$data[37,3] = 36584
$data[37,5] = 'M'
$data[37,6] = 'mking@gmail.com'

$data[38,0] = 39
$data[38,1] = 145646
$data[38,2] = 'Nalugala, Venessa Chebukwa'
$data[38,3] = 37208
$data[38,5] = 'F'
$data[38,6] = 'nchebukwa@gmail.com'

$data[39,0] = 40
$data[39,1] = 145813
$data[39,2] = 'Ndirangu, Denise Nyambura'
$data[39,3] = 37598
$data[39,5] = 'F'
$data[39,6] = 'nnyambura1@gmail.com'

$data[40,0] = 41
$data[40,1] = 139149
$data[40,2] = 'Ngahu, David Gitonga'
$data[40,3] = 36747
$data[40,5] = 'M'
$data[40,6] = 'ngitonga@gmail.com'

$data[41,0] = 42
$data[41,1] = 141733
$data[41,2] = 'Ngari, Sifa Gathoni'
$data[41,3] = 37787
$data[41,5] = 'F'
$data[41,6] = 'ngathoni@gmail.com'

$data[42,0] = 43
$data[42,1] = 134321
$data[42,2] = 'Nguthiru, Edwin Ndiritu'
$data[42,3] = 36895
$data[42,5] = 'M'
$data[42,6] = 'nndiritu@gmail.com'

$data[43,0] = 44
$data[43,1] = 145354
$data[43,2] = 'Nyamosi, Edmond Omwega'
$data[43,3] = 36993
$data[43,5] = 'M'
$data[43,6] = 'nomwega@gmail.com'

$data[44,0] = 45
$data[44,1] = 145536
$data[44,2] = 'Nyang''or, Olive Menorah'
$data[44,3] = 37449
$data[44,5] = 'F'
$data[44,6] = 'nmenorah@gmail.com'

$data[45,0] = 46
$data[45,1] = 122993
$data[45,2] = 'Omal, Warren'
$data[45,3] = 37810
$data[45,5] = 'M'
$data[45,6] = 'owarren@gmail.com'

$data[46,0] = 47
$data[46,1] = 145182
$data[46,2] = 'Omondi, Emmanuel Neville'
$data[46,3] = 37941
$data[46,5] = 'M'
$data[46,6] = 'oneville@gmail.com'

$data[47,0] = 48
$data[47,1] = 146533
$data[47,2] = 'Omondi, Winfred Achieng'
$data[47,3] = 36695
$data[47,5] = 'F'
$data[47,6] = 'oachieng@gmail.com'

$data[48,0] = 49
$data[48,1] = 131778
$data[48,2] = 'Otao, Davis Mokora'
$data[48,3] = 37046
$data[48,5] = 'M'
$data[48,6] = 'omokora@gmail.com'

$data[49,0] = 50
$data[49,1] = 146202
$data[49,2] = 'Rintaugu, Mugambi Nteere'
$data[49,3] = 36674
$data[49,5] = 'M'
$data[49,6] = 'rnteere@gmail.com'

$data[50,0] = 51
$data[50,1] = 146413
$data[50,2] = 'Rotich, Mercy Chepngetich'
$data[50,3] = 37812
$data[50,5] = 'F'
$data[50,6] = 'rchepngetich@gmail.com'

$data[51,0] = 52
$data[51,1] = 146254
$data[51,2] = 'Sehmi, Singh Jeevan'
$data[51,3] = 36915
$data[51,5] = 'M'
$data[51,6] = 'sjeevan@gmail.com'

$data[52,0] = 53
$data[52,1] = 144914
$data[52,2] = 'Vasani, Aman Upinkumar'
$data[52,3] = 37616
$data[52,5] = 'M'
$data[52,6] = 'vupinkumar@gmail.com'

$data[53,0] = 54
$data[53,1] = 135361
$data[53,2] = 'Wafula, Gideon Simiyu'
$data[53,3] = 37970
$data[53,5] = 'M'
$data[53,6] = 'wsimiyu@gmail.com'

$data[54,0] = 55
$data[54,1] = 145770
$data[54,2] = 'Wahu, Bridget Makena'
$data[54,3] = 37356
$data[54,5] = 'F'
$data[54,6] = 'wmakena@gmail.com'

$data[55,0] = 56
$data[55,1] = 145369
$data[55,2] = 'Wango, Michael Mundati'
$data[55,3] = 37596
$data[55,5] = 'M'
$data[55,6] = 'wmundati@gmail.com'

$data[56,0] = 57
$data[56,1] = 137938
$data[56,2] = 'Wangombe, Martin Mwangi'
$data[56,3] = 37476
$data[56,5] = 'M'
$data[56,6] = 'wmwangi@gmail.com'

$data[57,0] = 58
$data[57,1] = 145838
$data[57,2] = 'Wanyonyi, Brian Newton'
$data[57,3] = 37304
$data[57,5] = 'M'
$data[57,6] = 'wnewton@gmail.com'

$data[58,0] = 59
$data[58,1] = 138616
$data[58,2] = 'Wasike, Nicole Jones Nekesa'
$data[58,3] = 37945
$data[58,5] = 'F'
$data[58,6] = 'wnekesa@gmail.com'

$endRow = $startRow + $rowCount - 1
$targetRange = $ws.Range("A" + $startRow + ":G" + $endRow)
$targetRange.Value = $data

# Apply the same DoB date format used elsewhere in column D
$ws.Range("D" + $startRow + ":D" + $endRow).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Column E (Email Address) stays blank for the new rows, matching existing rows
$ws.Range("E" + $startRow + ":E" + $endRow).Style = "Normal"

Write-Output "Appended rows 66-124 to Sheet1"
